$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Remove the two rows for account 004471893 (Paula, 59677.67) and
# 005701765 (F, 52880.71), which sit immediately below the header row.
$ws.Range("A5:A6").EntireRow.Delete()
